$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.04271373187048222
$ws.Range("C2").Value = 0.04071648406533734
$ws.Range("D2").Value = 0.1494219747398047
$ws.Range("E2").Value = 10.19245300693656
$ws.Range("G2").Value = 10.42530519761218

# Row 3
$ws.Range("B3").Value = 3.286832544864788
$ws.Range("C3").Value = 1.655778082260271
$ws.Range("D3").Value = 3.537761648806719
$ws.Range("E3").Value = 10.19245300693656
$ws.Range("G3").Value = 18.67282528286833

# Row 4
$ws.Range("B4").Value = 3.286832544864788
$ws.Range("C4").Value = 1.655778082260271
$ws.Range("D4").Value = 0.1494219747398047
$ws.Range("E4").Value = 0.4942365360607697
$ws.Range("G4").Value = 5.586269137925634

# Row 5
$ws.Range("B5").Value = 0.2917716402565462
$ws.Range("C5").Value = 0.002571899574220771
$ws.Range("D5").Value = 0.1494219747398047
$ws.Range("E5").Value = 0.4942365360607697
$ws.Range("G5").Value = 0.9380020506313413

# Row 6
$ws.Range("B6").Value = 0.2917716402565462
$ws.Range("C6").Value = 1.655778082260271
$ws.Range("D6").Value = 3.537761648806719
$ws.Range("E6").Value = 10.19245300693656
$ws.Range("G6").Value = 15.67776437826009

# Row 7
$ws.Range("B7").Value = 3.286832544864788
$ws.Range("C7").Value = 1.655778082260271
$ws.Range("D7").Value = 3.537761648806719
$ws.Range("E7").Value = 0.4942365360607697
$ws.Range("G7").Value = 8.974608811992548

$wb.Save()
